$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$link1 = "https://www.genomeweb.com/molecular-diagnostics/source-bioscience-acquires-cambridge-clinical-laboratories"
$link2 = "https://www.360dx.com/molecular-diagnostics/source-bioscience-acquires-cambridge-clinical-laboratories"
$title = "Source BioScience Acquires Cambridge Clinical Laboratories"
$keyword = "digital pathology"

# Row 43
$ws.Range("A43").Value = $link1
$ws.Range("B43").Value = $keyword
$ws.Range("C43").Value = $title

# Row 44
$ws.Range("A44").Value = $link2
$ws.Range("B44").Value = $keyword
$ws.Range("C44").Value = $title

# Add hyperlinks (mirrors existing rows which carry rId hyperlink relationships)
$ws.Hyperlinks.Add($ws.Range("A43"), $link1)
$ws.Hyperlinks.Add($ws.Range("A44"), $link2)

# Match the style used by the other link cells in column A (Hyperlink cell style)
$ws.Range("A43").Style = "Hyperlink"
$ws.Range("A44").Style = "Hyperlink"
